$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D2:D51) to be treated as text so values
# like "1.00" or "327.47" are not coerced into numbers, matching the
# original inline-string cell contents.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '43.276.53'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '2.402.49'
$ws.Range('E3').Value = '  +5.62%  '
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').Value = '327.47'
$ws.Range('E5').Value = '  +7.88%  '
$ws.Range('D6').Value = '106.34'
$ws.Range('E6').Value = '  -6.52%  '
$ws.Range('D7').Value = '0.652'
$ws.Range('E7').Value = '  +2.85%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.655'
$ws.Range('E9').Value = '  +6.29%  '
$ws.Range('D10').Value = '42.28'
$ws.Range('E10').Value = '  -5.45%  '
$ws.Range('D11').Value = '0.0942'
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').Value = '8.77'
$ws.Range('E12').Value = '  -2.17%  '
$ws.Range('D13').Value = '1.05'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('D14').Value = '17.16'
$ws.Range('E14').Value = '  +10.46%  '
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = '2.764.63'
$ws.Range('E16').Value = '  +5.78%  '
$ws.Range('D17').Value = '2.393.61'
$ws.Range('E17').Value = '  +5.14%  '
$ws.Range('D18').Value = '43.268.84'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('D19').Value = '7.78'
$ws.Range('E19').Value = '  +7.42%  '
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('D21').Value = '77.15'
$ws.Range('E21').Value = '  +2.55%  '
$ws.Range('D22').Value = '3.71'
$ws.Range('E22').Value = '  +4.10%  '
$ws.Range('D23').Value = '274.19'
$ws.Range('E23').Value = '  +6.25%  '
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').Value = '9.69'
$ws.Range('E25').Value = '  +6.75%  '
$ws.Range('D26').Value = '11.84'
$ws.Range('E26').Value = '  +1.13%  '
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').Value = '23.28'
$ws.Range('E28').Value = '  +4.13%  '
$ws.Range('D29').Value = '176.65'
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('E30').Value = '  -1.98%  '
$ws.Range('D31').Value = '37.36'
$ws.Range('E31').Value = '  -2.41%  '
$ws.Range('D32').Value = '0.0942'
$ws.Range('E32').Value = '  +4.93%  '
$ws.Range('E33').Value = '  -0.55%  '
$ws.Range('D34').Value = '5.97'
$ws.Range('E35').Value = '  +4.94%  '
$ws.Range('D36').Value = '4.91'
$ws.Range('E36').Value = '  -4.17%  '
$ws.Range('E37').Value = '  -3.63%  '
$ws.Range('D38').Value = '0.0366'
$ws.Range('E38').Value = '  -3.51%  '
$ws.Range('D39').Value = '0.108'
$ws.Range('E39').Value = '  +2.82%  '
$ws.Range('D40').Value = '2.85'
$ws.Range('E40').Value = '  +15.59%  '
$ws.Range('E41').Value = '  +18.01%  '
$ws.Range('D42').Value = '0.237'
$ws.Range('E42').Value = '  +0.94%  '
$ws.Range('D43').Value = '70.42'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('D44').Value = '122.80'
$ws.Range('E44').Value = '  +13.80%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '92.20'
$ws.Range('E46').Value = '  +42.14%  '
$ws.Range('D47').Value = '12.38'
$ws.Range('E47').Value = '  -2.37%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').Value = '5.58'
$ws.Range('E48').Value = '  -0.72%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '9.35'
$ws.Range('E49').Value = '  +7.21%  '
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '0.493'
$ws.Range('E51').Value = '  +12.10%  '
